# Sprint 1 backlog: update the "Completed" (column D) burn-down figures
# for sprint 1 so the chart/summary reflect actual completed points.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Tasks that are no longer counted as complete this sprint -> clear D value
$ws.Range("D4").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("D21").ClearContents()

# Tasks with updated partial-completion point values
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 2

# Recalculate so the D28 SUM (and anything depending on it) is fresh
$excel.Calculate()

# Restore the selection to match the author's final cursor position
$ws.Range("D4").Select()
